$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.604.26'
$ws.Range('E2').Value = '  +2.00%  '
$ws.Range('D3').Value = '2.195.96'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '257.13'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '84.16'
$ws.Range('E6').Value = '  +11.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.619'
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.45'
$ws.Range('E10').Value = '  +7.94%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0916'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.22'
$ws.Range('E12').Value = '  +5.00%  '
$ws.Range('E13').Value = '  +1.97%  '
$ws.Range('D14').Value = '2.528.90'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.33'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '2.193.39'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.780'
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').Value = '43.501.98'
$ws.Range('E18').Value = '  +2.01%  '
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.68'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.89'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.35'
$ws.Range('E22').Value = '  +6.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '230.77'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.88'
$ws.Range('E24').Value = '  -5.54%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'WEMIXToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.50'
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.59'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('E28').Value = '  +3.05%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.72'
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '173.12'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.33'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.29'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.122'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.109'
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0357'
$ws.Range('E37').Value = '  +3.30%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.45'
$ws.Range('E38').Value = '  +3.42%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '12.32'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('E40').Value = '  +4.32%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.08'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '62.91'
$ws.Range('E42').Value = '  +4.66%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.44'
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.197'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '99.87'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0977'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.30'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.17'
$ws.Range('E48').Value = '  +3.66%  '
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.433'
$ws.Range('E50').Value = '  -6.31%  '
$ws.Range('E51').Value = '  +3.56%  '
